$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 4 odds (values refreshed in the source feed) ---
$ws.Range("G4").Value = 2.18
$ws.Range("I4").Value = 2.82
$ws.Range("J4").Value = 2.62
$ws.Range("K4").Value = 2.42
$ws.Range("L4").Value = 3.15
$ws.Range("M4").Value = 1.03
$ws.Range("P4").Value = 5.1
$ws.Range("Q4").Value = 1.42
$ws.Range("R4").Value = 2.67
$ws.Range("S4").Value = 1.24
$ws.Range("W4").Value = 14
$ws.Range("X4").Value = 15.5
$ws.Range("Z4").Value = 25
$ws.Range("AA4").Value = 15
$ws.Range("AD4").Value = 8.25
$ws.Range("AH4").Value = 16.5
$ws.Range("AI4").Value = 21
$ws.Range("AJ4").Value = 10.75
$ws.Range("AK4").Value = 37
$ws.Range("AL4").Value = 19.5
$ws.Range("AN4").Value = 4.8
$ws.Range("AO4").Value = 10.75
$ws.Range("AP4").Value = 14
$ws.Range("AW4").Value = 5.5
$ws.Range("AX4").Value = 14
$ws.Range("AZ4").Value = 50

# --- Append a new match (row 6) ---
$ws.Range("A6").Value = "r9LRh2ph"
$ws.Range("B6").Value = "31/10/2024"
$ws.Range("C6").Value = "12:10"
$ws.Range("D6").Value = "SAUDI ARABIA - SAUDI PROFESSIONAL LEAGUE"
$ws.Range("E6").Value = "Al Okhdood"
$ws.Range("F6").Value = "Al Orubah"
$ws.Range("G6").Value = 1.8
$ws.Range("H6").Value = 3.7
$ws.Range("I6").Value = 4
$ws.Range("J6").Value = 2.38
$ws.Range("K6").Value = 2.2
$ws.Range("L6").Value = 4.33
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 8.5
$ws.Range("O6").Value = 1.25
$ws.Range("P6").Value = 3.75
$ws.Range("Q6").Value = 1.83
$ws.Range("R6").Value = 1.98
$ws.Range("S6").Value = 1.36
$ws.Range("T6").Value = 3
$ws.Range("U6").Value = 1.8
$ws.Range("V6").Value = 1.91
$ws.Range("W6").Value = 7.5
$ws.Range("X6").Value = 9
$ws.Range("Y6").Value = 9
$ws.Range("Z6").Value = 15
$ws.Range("AA6").Value = 15
$ws.Range("AB6").Value = 26
$ws.Range("AC6").Value = 11
$ws.Range("AD6").Value = 7
$ws.Range("AE6").Value = 15
$ws.Range("AF6").Value = 51
$ws.Range("AG6").Value = 600
$ws.Range("AH6").Value = 12
$ws.Range("AI6").Value = 21
$ws.Range("AJ6").Value = 15
$ws.Range("AK6").Value = 41
$ws.Range("AL6").Value = 34
$ws.Range("AM6").Value = 41
$ws.Range("AN6").Value = 4
$ws.Range("AO6").Value = 9.5
$ws.Range("AP6").Value = 21
$ws.Range("AQ6").Value = 29
$ws.Range("AR6").Value = 51
$ws.Range("AS6").Value = 126
$ws.Range("AT6").Value = 3
$ws.Range("AU6").Value = 8
$ws.Range("AV6").Value = 51
$ws.Range("AW6").Value = 6
$ws.Range("AX6").Value = 21
$ws.Range("AY6").Value = 29
$ws.Range("AZ6").Value = 67
$ws.Range("BA6").Value = 81
$ws.Range("BB6").Value = 151
$ws.Range("BC6").Value = 81
$ws.Range("BD6").Value = 81
